$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '50.754.65'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -1.36%  '
$ws.Range('E2').ClearFormats()
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.921.50'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -2.27%  '
$ws.Range('E3').ClearFormats()
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('E4').ClearFormats()
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '373.18'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -2.96%  '
$ws.Range('E5').ClearFormats()
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '99.81'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -3.33%  '
$ws.Range('E6').ClearFormats()
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -1.52%  '
$ws.Range('E7').ClearFormats()
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E8').ClearFormats()
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.583'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -1.67%  '
$ws.Range('E9').ClearFormats()
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '35.83'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -3.46%  '
$ws.Range('E10').ClearFormats()
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -0.68%  '
$ws.Range('E11').ClearFormats()
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -0.26%  '
$ws.Range('E12').ClearFormats()
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.387.64'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -1.96%  '
$ws.Range('E13').ClearFormats()
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '17.90'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -2.24%  '
$ws.Range('E14').ClearFormats()
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.52'
$ws.Range('D15').ClearFormats()
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.917.10'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -2.20%  '
$ws.Range('E16').ClearFormats()
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '11.05'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +48.98%  '
$ws.Range('E17').ClearFormats()
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.987'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -2.03%  '
$ws.Range('E18').ClearFormats()
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '50.709.63'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -1.35%  '
$ws.Range('E19').ClearFormats()
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.03'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -7.12%  '
$ws.Range('E20').ClearFormats()
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.34'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -4.08%  '
$ws.Range('E21').ClearFormats()
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.0₃0949'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -1.13%  '
$ws.Range('E22').ClearFormats()
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '68.84'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.43%  '
$ws.Range('E23').ClearFormats()
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '264.37'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +0.98%  '
$ws.Range('E24').ClearFormats()
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +7.11%  '
$ws.Range('E25').ClearFormats()
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.01'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -2.42%  '
$ws.Range('E26').ClearFormats()
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.36'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -3.77%  '
$ws.Range('E27').ClearFormats()
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '25.43'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -2.20%  '
$ws.Range('E29').ClearFormats()
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -5.02%  '
$ws.Range('E30').ClearFormats()
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.108'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -9.18%  '
$ws.Range('E31').ClearFormats()
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '9.95'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +0.56%  '
$ws.Range('E32').ClearFormats()
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '50.76'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -0.43%  '
$ws.Range('E33').ClearFormats()
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '32.99'
$ws.Range('D35').ClearFormats()
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0436'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -4.35%  '
$ws.Range('E36').ClearFormats()
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.10'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +3.48%  '
$ws.Range('E38').ClearFormats()
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -1.09%  '
$ws.Range('E39').ClearFormats()
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '16.36'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -3.86%  '
$ws.Range('E40').ClearFormats()
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -2.21%  '
$ws.Range('E41').ClearFormats()
$ws.Range('B42').NumberFormat = "@"
$ws.Range('B42').Value = 'Monero'
$ws.Range('B42').ClearFormats()
$ws.Range('C42').NumberFormat = "@"
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('C42').ClearFormats()
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '123.73'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +0.57%  '
$ws.Range('E42').ClearFormats()
$ws.Range('B43').NumberFormat = "@"
$ws.Range('B43').Value = 'Stacks'
$ws.Range('B43').ClearFormats()
$ws.Range('C43').NumberFormat = "@"
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('C43').ClearFormats()
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.43'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -6.06%  '
$ws.Range('E43').ClearFormats()
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '20.88'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -3.39%  '
$ws.Range('E44').ClearFormats()
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -0.90%  '
$ws.Range('E45').ClearFormats()
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.35'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +1.51%  '
$ws.Range('E46').ClearFormats()
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -1.34%  '
$ws.Range('E47').ClearFormats()
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.267'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -2.36%  '
$ws.Range('E48').ClearFormats()
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.979.30'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -2.81%  '
$ws.Range('E49').ClearFormats()
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0322'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -3.34%  '
$ws.Range('E50').ClearFormats()
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '5.17'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +0.61%  '
$ws.Range('E51').ClearFormats()
